$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so that numeric-looking
# strings (e.g. "10.50", "275.00") keep their exact formatting instead of
# being auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "54.792.20"
$ws.Range("E2").Value = "  +5.77%  "
$ws.Range("D3").Value = "3.190.47"
$ws.Range("E3").Value = "  +3.15%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "399.11"
$ws.Range("E5").Value = "  +2.49%  "
$ws.Range("D6").Value = "109.97"
$ws.Range("E6").Value = "  +5.73%  "
$ws.Range("D7").Value = "0.551"
$ws.Range("E7").Value = "  +0.94%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "0.619"
$ws.Range("E9").Value = "  +4.96%  "
$ws.Range("D10").Value = "39.12"
$ws.Range("E10").Value = "  +5.21%  "
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("D12").Value = "0.0887"
$ws.Range("E12").Value = "  +2.39%  "
$ws.Range("D13").Value = "3.681.51"
$ws.Range("E13").Value = "  +3.02%  "
$ws.Range("D14").Value = "19.24"
$ws.Range("E14").Value = "  +2.35%  "
$ws.Range("D15").Value = "8.07"
$ws.Range("E15").Value = "  +2.63%  "
$ws.Range("E16").Value = "  +8.08%  "
$ws.Range("D17").Value = "3.183.42"
$ws.Range("D18").Value = "10.50"
$ws.Range("E18").Value = "  -2.08%  "
$ws.Range("D19").Value = "54.547.54"
$ws.Range("E19").Value = "  +5.16%  "
$ws.Range("D20").Value = "3.31"
$ws.Range("E20").Value = "  +4.18%  "
$ws.Range("D21").Value = "12.93"
$ws.Range("E21").Value = "  +3.05%  "
$ws.Range("D22").Value = "0.0₃0999"
$ws.Range("E22").Value = "  +2.72%  "
$ws.Range("D23").Value = "71.64"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("D24").Value = "275.00"
$ws.Range("E24").Value = "  +2.04%  "
$ws.Range("E25").Value = "  +3.00%  "
$ws.Range("D26").Value = "8.00"
$ws.Range("E26").Value = "  -2.57%  "
$ws.Range("D27").Value = "27.78"
$ws.Range("E27").Value = "  +2.68%  "
$ws.Range("D28").Value = "7.36"
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("D29").Value = "0.170"
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "0.113"
$ws.Range("E31").Value = "  +4.67%  "
$ws.Range("E32").Value = "  +6.55%  "
$ws.Range("E33").Value = "  +10.41%  "
$ws.Range("D34").Value = "37.04"
$ws.Range("E34").Value = "  +3.68%  "
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("D36").Value = "50.68"
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("D37").Value = "3.66"
$ws.Range("E37").Value = "  +8.17%  "
$ws.Range("D38").Value = "0.997"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").Value = "2.88"
$ws.Range("E39").Value = "  +11.21%  "
$ws.Range("D40").Value = "4.11"
$ws.Range("E40").Value = "  +9.63%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "1.92"
$ws.Range("E41").Value = "  +1.93%  "
$ws.Range("D42").Value = "0.291"
$ws.Range("E42").Value = "  -2.04%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").Value = "17.34"
$ws.Range("E43").Value = "  +2.01%  "
$ws.Range("D44").Value = "130.20"
$ws.Range("E44").Value = "  +2.63%  "
$ws.Range("D45").Value = "0.118"
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("D46").Value = "22.30"
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("D47").Value = "2.45"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("D49").Value = "2.093.93"
$ws.Range("E49").Value = "  +2.69%  "
$ws.Range("D50").Value = "0.0345"
$ws.Range("E50").Value = "  +7.00%  "
$ws.Range("D51").Value = "0.0502"
$ws.Range("E51").Value = "  +10.26%  "
